# Auto-generated PowerShell Excel COM-interop script
# Updates the "想去人数" (interested-count) column F values across all 4 sheets
# to match the regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F10").Value = 4834
$ws.Range("F11").Value = 6216
$ws.Range("F12").Value = 862
$ws.Range("F13").Value = 80
$ws.Range("F14").Value = 1375
$ws.Range("F17").Value = 6746
$ws.Range("F21").Value = 4516
$ws.Range("F22").Value = 356
$ws.Range("F24").Value = 744
$ws.Range("F28").Value = 1119
$ws.Range("F29").Value = 157
$ws.Range("F33").Value = 354
$ws.Range("F35").Value = 1947
$ws.Range("F36").Value = 186
$ws.Range("F39").Value = 1295
$ws.Range("F42").Value = 52
$ws.Range("F43").Value = 1031
$ws.Range("F44").Value = 1320
$ws.Range("F45").Value = 35

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 443
$ws.Range("F11").Value = 2
$ws.Range("F13").Value = 232
$ws.Range("F26").Value = 242
$ws.Range("F27").Value = 11

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F11").Value = 1667
$ws.Range("F12").Value = 2000
$ws.Range("F13").Value = 459

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 1667
$ws.Range("F10").Value = 2000
$ws.Range("F11").Value = 4834
$ws.Range("F12").Value = 443
$ws.Range("F14").Value = 862
$ws.Range("F15").Value = 80
$ws.Range("F17").Value = 1375
$ws.Range("F20").Value = 6746
$ws.Range("F24").Value = 4516
$ws.Range("F25").Value = 356
$ws.Range("F26").Value = 744
$ws.Range("F28").Value = 1119
$ws.Range("F29").Value = 157
$ws.Range("F34").Value = 354
$ws.Range("F36").Value = 1947
$ws.Range("F37").Value = 186
$ws.Range("F41").Value = 1295
$ws.Range("F46").Value = 1320
